$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.109.39"
$ws.Range("E2").Value = "  +0.08%  "

$ws.Range("D3").Value = "'1.830.72"
$ws.Range("E3").Value = "  -0.39%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'243.08"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("D6").Value = "'0.6252"
$ws.Range("E6").Value = "  -0.53%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").Value = "'0.07472"
$ws.Range("E8").Value = "  -1.53%  "

$ws.Range("D9").Value = "'0.2925"
$ws.Range("E9").Value = "  -0.22%  "

$ws.Range("D10").Value = "'23.29"
$ws.Range("E10").Value = "  +2.98%  "

$ws.Range("E11").Value = "  -0.65%  "

$ws.Range("D12").Value = "'1.845.13"
$ws.Range("E12").Value = "  +0.31%  "

$ws.Range("D13").Value = "'5.011"
$ws.Range("E13").Value = "  +0.75%  "

$ws.Range("D14").Value = "'0.6671"
$ws.Range("E14").Value = "  +0.09%  "

$ws.Range("D15").Value = "'82.60"
$ws.Range("E15").Value = "  -0.73%  "

$ws.Range("D16").Value = "'0.000009360"
$ws.Range("E16").Value = "  -6.59%  "

$ws.Range("D17").Value = "'5.970"
$ws.Range("E17").Value = "  -1.61%  "

$ws.Range("D18").Value = "'29.095.89"
$ws.Range("E18").Value = "  +0.02%  "

$ws.Range("D19").Value = "'2.068.47"
$ws.Range("E19").Value = "  -0.88%  "

$ws.Range("D20").Value = "'12.59"
$ws.Range("E20").Value = "  +1.29%  "

$ws.Range("D21").Value = "'222.90"
$ws.Range("E21").Value = "  -1.85%  "

$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("E23").Value = "  -1.42%  "

$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").Value = "'1.002"
$ws.Range("E24").Value = "  -0.28%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'160.23"
$ws.Range("E25").Value = "  +0.35%  "

$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").Value = "'0.1394"
$ws.Range("E26").Value = "  +0.55%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'8.499"
$ws.Range("E27").Value = "  -0.17%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'17.88"
$ws.Range("E28").Value = "  -0.31%  "

$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'1.491"
$ws.Range("E29").Value = "  -0.32%  "

$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.05835"
$ws.Range("E30").Value = "  +11.08%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.157"
$ws.Range("E31").Value = "  +1.33%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'4.115"
$ws.Range("E32").Value = "  +2.40%  "

$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").Value = "'1.213"
$ws.Range("E33").Value = "  +1.45%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7398"
$ws.Range("E34").Value = "  +0.35%  "

$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'1.828"
$ws.Range("E35").Value = "  -1.04%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.137"
$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.673"
$ws.Range("E37").Value = "  -0.29%  "

$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "'1.228.14"
$ws.Range("E38").Value = "  -1.52%  "

$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.765"
$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01769"
$ws.Range("E40").Value = "  -0.85%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'6.493"
$ws.Range("E41").Value = "  +1.98%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.8942"
$ws.Range("E42").Value = "  -0.74%  "

$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "'1.001"
$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'102.14"
$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'66.02"
$ws.Range("E45").Value = "  +2.73%  "

$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "'1.962.38"
$ws.Range("E46").Value = "  -1.19%  "

$ws.Range("D47").Value = "'0.00000000123"
$ws.Range("E47").Value = "  +0.11%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.5093"
$ws.Range("E48").Value = "  -0.60%  "

$ws.Range("B49").Value = "XinFinNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D49").Value = "'0.07535"
$ws.Range("E49").Value = "  +12.95%  "

$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").Value = "'0.4063"
$ws.Range("E50").Value = "  +0.51%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'8.997"
$ws.Range("E51").Value = "  +1.43%  "

